$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(82, 8).Value = 1557.875
$ws.Cells.Item(82, 9).Value = 1579.1666
$ws.Cells.Item(82, 11).Value = 4737.4998
$ws.Cells.Item(82, 13).Value = -4331.4998
$ws.Cells.Item(85, 8).Value = 1557.875
$ws.Cells.Item(85, 9).Value = 1579.1666
$ws.Cells.Item(85, 11).Value = 4737.4998
$ws.Cells.Item(85, 13).Value = -3333.4998
$ws.Cells.Item(99, 8).Value = 1766.3334
$ws.Cells.Item(99, 9).Value = 199
$ws.Cells.Item(99, 10).Value = 2079.8
$ws.Cells.Item(99, 11).Value = 597
$ws.Cells.Item(99, 12).Value = 6239.400000000001
$ws.Cells.Item(99, 13).Value = 901
$ws.Cells.Item(99, 14).Value = -9235.400000000001
$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 14).ClearContents()
$ws.Cells.Item(115, 8).Value = 846.1
$ws.Cells.Item(115, 9).Value = 220.6
$ws.Cells.Item(115, 11).Value = 661.8
$ws.Cells.Item(115, 13).Value = 905.2
$ws.Cells.Item(137, 8).Value = 3245.0417
$ws.Cells.Item(137, 10).Value = 4344.846
$ws.Cells.Item(137, 12).Value = 13034.538
$ws.Cells.Item(137, 14).Value = -18134.538
$ws.Cells.Item(138, 8).Value = 3736.0322
$ws.Cells.Item(138, 9).Value = 1545.7709
$ws.Cells.Item(138, 10).Value = 6072.311
$ws.Cells.Item(138, 11).Value = 4637.3127
$ws.Cells.Item(138, 12).Value = 18216.933
$ws.Cells.Item(138, 13).Value = 502.6872999999996
$ws.Cells.Item(138, 14).Value = -28496.933
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4124.512
$ws.Cells.Item(32, 9).Value = 3308.639
$ws.Cells.Item(32, 11).Value = 3308.639
$ws.Cells.Item(32, 13).Value = -3021.639
$ws.Cells.Item(74, 8).Value = 2670.5806
$ws.Cells.Item(74, 9).Value = 2531.5217
$ws.Cells.Item(74, 10).Value = 3070.375
$ws.Cells.Item(74, 11).Value = 2531.5217
$ws.Cells.Item(74, 12).Value = 3070.375
$ws.Cells.Item(74, 13).Value = -1657.5217
$ws.Cells.Item(74, 14).Value = -4818.375
$ws.Cells.Item(77, 8).Value = 2670.5806
$ws.Cells.Item(77, 9).Value = 2531.5217
$ws.Cells.Item(77, 10).Value = 3070.375
$ws.Cells.Item(77, 11).Value = 12657.6085
$ws.Cells.Item(77, 12).Value = 15351.875
$ws.Cells.Item(77, 13).Value = -8289.608499999998
$ws.Cells.Item(77, 14).Value = -24087.875
$ws.Cells.Item(97, 8).Value = 1420.6774
$ws.Cells.Item(97, 9).Value = 1552.0385
$ws.Cells.Item(97, 10).Value = 737.6
$ws.Cells.Item(97, 11).Value = 1552.0385
$ws.Cells.Item(97, 12).Value = 737.6
$ws.Cells.Item(97, 13).Value = -1056.0385
$ws.Cells.Item(97, 14).Value = -1729.6
$ws.Cells.Item(110, 8).Value = 835343
$ws.Cells.Item(110, 9).Value = 1002011.6
$ws.Cells.Item(110, 10).Value = 2000
$ws.Cells.Item(110, 11).Value = 1002011.6
$ws.Cells.Item(110, 12).Value = 2000
$ws.Cells.Item(110, 13).Value = -999966.6
$ws.Cells.Item(110, 14).Value = -6090
$ws.Cells.Item(122, 8).Value = 4742.091
$ws.Cells.Item(122, 9).Value = 3947.5
$ws.Cells.Item(122, 10).Value = 4918.6665
$ws.Cells.Item(122, 11).Value = 11842.5
$ws.Cells.Item(122, 12).Value = 14755.9995
$ws.Cells.Item(122, 13).Value = -9392.5
$ws.Cells.Item(122, 14).Value = -19655.9995
$ws.Cells.Item(132, 8).Value = 6349.8076
$ws.Cells.Item(132, 9).Value = 2182.0588
$ws.Cells.Item(132, 10).Value = 14222.223
$ws.Cells.Item(132, 11).Value = 6546.176399999999
$ws.Cells.Item(132, 12).Value = 42666.669
$ws.Cells.Item(132, 13).Value = -4016.176399999999
$ws.Cells.Item(132, 14).Value = -47726.669
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3321.4
$ws.Cells.Item(20, 9).Value = 2502.4167
$ws.Cells.Item(20, 11).Value = 2502.4167
$ws.Cells.Item(20, 13).Value = -2255.4167
$ws.Cells.Item(134, 8).Value = 23778.56
$ws.Cells.Item(134, 9).Value = 2673.919
$ws.Cells.Item(134, 10).Value = 83845.62
$ws.Cells.Item(134, 11).Value = 8021.757
$ws.Cells.Item(134, 12).Value = 251536.86
$ws.Cells.Item(134, 13).Value = -5486.757
$ws.Cells.Item(134, 14).Value = -256606.86
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 685.1875
$ws.Cells.Item(7, 9).Value = 688
$ws.Cells.Item(7, 11).Value = 688
$ws.Cells.Item(7, 13).Value = -575
$ws.Cells.Item(125, 8).Value = 55163
$ws.Cells.Item(125, 10).Value = 60326
$ws.Cells.Item(125, 12).Value = 60326
$ws.Cells.Item(125, 14).Value = -65246
$ws.Cells.Item(132, 8).Value = 3506.805
$ws.Cells.Item(132, 9).Value = 2601.0322
$ws.Cells.Item(132, 11).Value = 7803.096600000001
$ws.Cells.Item(132, 13).Value = -5273.096600000001
$ws.Cells.Item(134, 8).Value = 221219.77
$ws.Cells.Item(134, 9).Value = 2802.6667
$ws.Cells.Item(134, 10).Value = 531601.9399999999
$ws.Cells.Item(134, 11).Value = 8408.000100000001
$ws.Cells.Item(134, 12).Value = 1594805.82
$ws.Cells.Item(134, 13).Value = -5873.000100000001
$ws.Cells.Item(134, 14).Value = -1599875.82
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 2034811.9
$ws.Cells.Item(4, 10).Value = 46249.168
$ws.Cells.Item(4, 12).Value = 138747.504
$ws.Cells.Item(4, 14).Value = -138971.504
$ws.Cells.Item(36, 8).Value = 1817.1666
$ws.Cells.Item(36, 9).Value = 1475
$ws.Cells.Item(36, 10).Value = 2501.5
$ws.Cells.Item(36, 11).Value = 4425
$ws.Cells.Item(36, 12).Value = 7504.5
$ws.Cells.Item(36, 13).Value = -4256
$ws.Cells.Item(36, 14).Value = -7842.5
$ws.Cells.Item(41, 8).Value = 334.25
$ws.Cells.Item(41, 9).Value = 334.25
$ws.Cells.Item(41, 11).Value = 1002.75
$ws.Cells.Item(41, 13).Value = -664.75
$ws.Cells.Item(80, 8).Value = 7166.6665
$ws.Cells.Item(80, 10).Value = 7500
$ws.Cells.Item(80, 12).Value = 22500
$ws.Cells.Item(80, 14).Value = -24372
$ws.Cells.Item(83, 8).Value = 7166.6665
$ws.Cells.Item(83, 10).Value = 7500
$ws.Cells.Item(83, 12).Value = 67500
$ws.Cells.Item(83, 14).Value = -76860
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(7, 8).Value = 1549999.1
$ws.Cells.Item(7, 10).Value = 1549999.1
$ws.Cells.Item(7, 12).Value = 1549999.1
$ws.Cells.Item(7, 14).Value = -1550223.1
$ws.Cells.Item(8, 8).Value = 1549999.1
$ws.Cells.Item(8, 10).Value = 1549999.1
$ws.Cells.Item(8, 12).Value = 1549999.1
$ws.Cells.Item(8, 14).Value = -1550277.1
$ws.Cells.Item(80, 8).Value = 1180564.1
$ws.Cells.Item(80, 9).Value = 913235
$ws.Cells.Item(80, 11).Value = 913235
$ws.Cells.Item(80, 13).Value = -912237
$ws.Cells.Item(83, 8).Value = 1180564.1
$ws.Cells.Item(83, 9).Value = 913235
$ws.Cells.Item(83, 11).Value = 4566175
$ws.Cells.Item(83, 13).Value = -4561183
$ws.Cells.Item(102, 8).Value = 1693.9546
$ws.Cells.Item(102, 9).Value = 1264.5625
$ws.Cells.Item(102, 10).Value = 2839
$ws.Cells.Item(102, 11).Value = 1264.5625
$ws.Cells.Item(102, 12).Value = 2839
$ws.Cells.Item(102, 13).Value = 357.4375
$ws.Cells.Item(102, 14).Value = -6083
$ws.Cells.Item(122, 8).Value = 3275.3076
$ws.Cells.Item(122, 9).Value = 3275.3076
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 9825.9228
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -7375.9228
$ws.Cells.Item(122, 14).ClearContents()
$ws.Cells.Item(123, 8).Value = 0
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 12).Value = 0
$ws.Cells.Item(123, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 6096.75
$ws.Cells.Item(126, 9).Value = 4233.4287
$ws.Cells.Item(126, 10).Value = 8705.4
$ws.Cells.Item(126, 11).Value = 12700.2861
$ws.Cells.Item(126, 12).Value = 26116.2
$ws.Cells.Item(126, 13).Value = -10230.2861
$ws.Cells.Item(126, 14).Value = -31056.2
$ws.Cells.Item(132, 8).Value = 738375.2
$ws.Cells.Item(132, 9).Value = 1254957
$ws.Cells.Item(132, 11).Value = 3764871
$ws.Cells.Item(132, 13).Value = -3762341
$ws.Cells.Item(134, 8).Value = 49997.75
$ws.Cells.Item(134, 10).Value = 49997.75
$ws.Cells.Item(134, 12).Value = 149993.25
$ws.Cells.Item(134, 14).Value = -155063.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 675.8333
$ws.Cells.Item(16, 9).Value = 663.75
$ws.Cells.Item(16, 11).Value = 663.75
$ws.Cells.Item(16, 13).Value = -493.75
$ws.Cells.Item(40, 8).Value = 336001.34
$ws.Cells.Item(40, 9).Value = 336001.34
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 336001.34
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = -335865.34
$ws.Cells.Item(40, 14).ClearContents()
$ws.Cells.Item(46, 8).Value = 4657.048
$ws.Cells.Item(46, 10).Value = 5089.9
$ws.Cells.Item(46, 12).Value = 5089.9
$ws.Cells.Item(46, 14).Value = -5465.9
$ws.Cells.Item(102, 8).Value = 46853.332
$ws.Cells.Item(102, 10).Value = 46853.332
$ws.Cells.Item(102, 12).Value = 46853.332
$ws.Cells.Item(102, 14).Value = -53343.332
$ws.Cells.Item(132, 8).Value = 3975.15
$ws.Cells.Item(132, 9).Value = 3455.0312
$ws.Cells.Item(132, 11).Value = 10365.0936
$ws.Cells.Item(132, 13).Value = -7835.0936
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(9, 8).Value = 3074
$ws.Cells.Item(9, 9).Value = 2111
$ws.Cells.Item(9, 11).Value = 2111
$ws.Cells.Item(9, 13).Value = -1971
$ws.Cells.Item(12, 8).Value = 10006
$ws.Cells.Item(12, 9).Value = 10006
$ws.Cells.Item(12, 11).Value = 10006
$ws.Cells.Item(12, 13).Value = -9864
$ws.Cells.Item(14, 8).Value = 20798.6
$ws.Cells.Item(14, 9).Value = 20748.25
$ws.Cells.Item(14, 11).Value = 20748.25
$ws.Cells.Item(14, 13).Value = -20580.25
$ws.Cells.Item(62, 8).Value = 7999.8335
$ws.Cells.Item(62, 9).Value = 5599.8
$ws.Cells.Item(62, 10).Value = 9714.143
$ws.Cells.Item(62, 11).Value = 5599.8
$ws.Cells.Item(62, 12).Value = 9714.143
$ws.Cells.Item(62, 13).Value = -4975.8
$ws.Cells.Item(62, 14).Value = -10962.143
$ws.Cells.Item(65, 8).Value = 7999.8335
$ws.Cells.Item(65, 9).Value = 5599.8
$ws.Cells.Item(65, 10).Value = 9714.143
$ws.Cells.Item(65, 11).Value = 27999
$ws.Cells.Item(65, 12).Value = 48570.715
$ws.Cells.Item(65, 13).Value = -24879
$ws.Cells.Item(65, 14).Value = -54810.715
$ws.Cells.Item(122, 8).Value = 62505612
$ws.Cells.Item(122, 9).Value = 166670610
$ws.Cells.Item(122, 11).Value = 500011830
$ws.Cells.Item(122, 13).Value = -500009380
